# feat: ease opening waves and refresh skill visuals
#
# Rebalances the early "abyssal nightfall" wave table (lower opening
# difficulty / timestamps) and refreshes enemy formations + flavor text,
# splitting what used to be a single row 11 ("40060005"/line) into a
# revised row 11 (new "40060004"/cross data) plus a brand-new row 12
# that carries the old line-formation entry forward under serial 0007.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data columns (A:J) in this sheet are persisted as text, including
# cells that merely look numeric ("60", "08", "40060007", "12.0", ...).
# Plain `.Value =` assignment would let Excel's type-inference turn those
# into real numbers, so every numeric-looking cell we touch gets its
# NumberFormat forced to "@" (text) right before the value is written.
# Cells we never touch are left completely alone.

function Set-TextValue($a1, $value) {
    $ws.Range($a1).NumberFormat = "@"
    $ws.Range($a1).Value = $value
}

# --- Row 6 (serial 0001): only a subset of cells changed ---
Set-TextValue "E6" "20"
Set-TextValue "F6" "40060007"
Set-TextValue "G6" "8"
Set-TextValue "H6" "12.0"
$ws.Range("J6").Value = "虚潮行者只会缠斗，练习走位与护盾。"

# --- Row 7 (serial 0002) ---
Set-TextValue "D7" "28"
Set-TextValue "E7" "22"
Set-TextValue "F7" "40060006"
Set-TextValue "H7" "8.5"
$ws.Range("I7").Value = "swarm"
$ws.Range("J7").Value = "碎影涌入，速度快但血量低。"

# --- Row 8 (serial 0003) ---
Set-TextValue "D8" "54"
Set-TextValue "E8" "26"
Set-TextValue "F8" "40060001"
Set-TextValue "G8" "12"
Set-TextValue "H8" "10.5"
$ws.Range("I8").Value = "ring"
$ws.Range("J8").Value = "餍爬者投掷腐质胆汁，学会绕开溅射区。"

# --- Row 9 (serial 0004) ---
Set-TextValue "D9" "86"
Set-TextValue "E9" "30"
Set-TextValue "F9" "40060002"
Set-TextValue "G9" "10"
Set-TextValue "H9" "9.0"
$ws.Range("I9").Value = "arc"
$ws.Range("J9").Value = "歌祭徒三连音袭击，持续走位或击杀增幅者。"

# --- Row 10 (serial 0005) ---
Set-TextValue "D10" "124"
Set-TextValue "E10" "32"
Set-TextValue "F10" "40060003"
Set-TextValue "G10" "8"
Set-TextValue "H10" "8.5"
$ws.Range("I10").Value = "cone"
$ws.Range("J10").Value = "呼嚎者扇形震荡波附理智流失，利用空挡反击。"

# --- Row 11 (serial 0006): now carries what used to be the row-9 enemy data ---
Set-TextValue "D11" "164"
Set-TextValue "E11" "34"
Set-TextValue "F11" "40060004"
Set-TextValue "G11" "6"
Set-TextValue "H11" "13.5"
$ws.Range("I11").Value = "cross"
$ws.Range("J11").Value = "虚壳哨兵有 0.6 秒警示后扫射光束，注意走位。"

# --- Row 12 (new, serial 0007): carries the former row-11 line/boss data ---
Set-TextValue "A12" "60"
Set-TextValue "B12" "08"
Set-TextValue "C12" "0007"
Set-TextValue "D12" "206"
Set-TextValue "E12" "40"
Set-TextValue "F12" "40060005"
Set-TextValue "G12" "8"
Set-TextValue "H12" "10.5"
$ws.Range("I12").Value = "line"
$ws.Range("J12").Value = "掘锚者冲撞灯塔并引发震波，及时打断或闪避。"

Write-Host "waves.xlsx: rebalanced rows 6-12 (dimension now A4:J12)"
